# Auto-generated edit script applying the Moogle_Profits market-data refresh
# (values sourced from the scheduled runner's latest market snapshot).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6927.154
$ws.Range("I76").Value = 5399.8
$ws.Range("K76").Value = 5399.8
$ws.Range("M76").Value = -5084.8
$ws.Range("H79").Value = 6927.154
$ws.Range("I79").Value = 5399.8
$ws.Range("K79").Value = 5399.8
$ws.Range("M79").Value = -4307.8
$ws.Range("H80").Value = 1276.6
$ws.Range("I80").Value = 405.75
$ws.Range("J80").Value = 1857.1666
$ws.Range("K80").Value = 1217.25
$ws.Range("L80").Value = 5571.4998
$ws.Range("M80").Value = -219.25
$ws.Range("N80").Value = -7567.4998
$ws.Range("H83").Value = 1276.6
$ws.Range("I83").Value = 405.75
$ws.Range("J83").Value = 1857.1666
$ws.Range("K83").Value = 3651.75
$ws.Range("L83").Value = 16714.4994
$ws.Range("M83").Value = 1340.25
$ws.Range("N83").Value = -26698.4994
$ws.Range("H96").Value = 288.4
$ws.Range("I96").Value = 315.5
$ws.Range("J96").Value = 180
$ws.Range("K96").Value = 946.5
$ws.Range("L96").Value = 540
$ws.Range("M96").Value = 426.5
$ws.Range("N96").Value = -3286
$ws.Range("H112").Value = 6750.3335
$ws.Range("I112").Value = 2387
$ws.Range("J112").Value = 7997
$ws.Range("K112").Value = 7161
$ws.Range("L112").Value = 23991
$ws.Range("M112").Value = -6053
$ws.Range("N112").Value = -26207

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4407.0835
$ws.Range("I61").Value = 4004.1785
$ws.Range("J61").Value = 5817.25
$ws.Range("K61").Value = 4004.1785
$ws.Range("L61").Value = 5817.25
$ws.Range("M61").Value = -3792.1785
$ws.Range("N61").Value = -6241.25
$ws.Range("H122").Value = 4398.0586
$ws.Range("I122").Value = 5568.4443
$ws.Range("K122").Value = 16705.3329
$ws.Range("M122").Value = -14255.3329
$ws.Range("H136").Value = 4407.0835
$ws.Range("I136").Value = 4004.1785
$ws.Range("J136").Value = 5817.25
$ws.Range("K136").Value = 12012.5355
$ws.Range("L136").Value = 17451.75
$ws.Range("M136").Value = -9462.5355
$ws.Range("N136").Value = -22551.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 200
$ws.Range("M17").Value = -26
$ws.Range("H31").Value = 8044.7676
$ws.Range("I31").Value = 4493.778
$ws.Range("K31").Value = 4493.778
$ws.Range("M31").Value = -4198.778
$ws.Range("H34").Value = 8044.7676
$ws.Range("I34").Value = 4493.778
$ws.Range("K34").Value = 4493.778
$ws.Range("M34").Value = -4291.778
$ws.Range("H39").Value = 14810.167
$ws.Range("I39").Value = 10965.5
$ws.Range("J39").Value = 22499.5
$ws.Range("K39").Value = 10965.5
$ws.Range("L39").Value = 22499.5
$ws.Range("M39").Value = -10574.5
$ws.Range("N39").Value = -23281.5
$ws.Range("H49").Value = 14810.167
$ws.Range("I49").Value = 10965.5
$ws.Range("J49").Value = 22499.5
$ws.Range("K49").Value = 10965.5
$ws.Range("L49").Value = 22499.5
$ws.Range("M49").Value = -10783.5
$ws.Range("N49").Value = -22863.5
$ws.Range("H62").Value = 6899.8423
$ws.Range("I62").Value = 6811
$ws.Range("J62").Value = 6931.5713
$ws.Range("K62").Value = 6811
$ws.Range("L62").Value = 6931.5713
$ws.Range("M62").Value = -6187
$ws.Range("N62").Value = -8179.5713
$ws.Range("H65").Value = 6899.8423
$ws.Range("I65").Value = 6811
$ws.Range("J65").Value = 6931.5713
$ws.Range("K65").Value = 34055
$ws.Range("L65").Value = 34657.85649999999
$ws.Range("M65").Value = -30935
$ws.Range("N65").Value = -40897.85649999999
$ws.Range("H134").Value = 2015.4849
$ws.Range("I134").Value = 1411.8334
$ws.Range("J134").Value = 3625.2222
$ws.Range("K134").Value = 4235.5002
$ws.Range("L134").Value = 10875.6666
$ws.Range("M134").Value = -1700.5002
$ws.Range("N134").Value = -15945.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 613.73334
$ws.Range("I8").Value = 613.73334
$ws.Range("K8").Value = 1841.20002
$ws.Range("M8").Value = -1702.20002
$ws.Range("H11").Value = 15118.777
$ws.Range("I11").Value = 125000
$ws.Range("K11").Value = 375000
$ws.Range("M11").Value = -374860
$ws.Range("H61").Value = 272.5
$ws.Range("I61").Value = 213.33333
$ws.Range("J61").Value = 450
$ws.Range("K61").Value = 639.99999
$ws.Range("L61").Value = 1350
$ws.Range("M61").Value = -424.99999
$ws.Range("N61").Value = -1780
$ws.Range("H109").Value = 1083.5
$ws.Range("I109").Value = 1130
$ws.Range("J109").Value = 944
$ws.Range("K109").Value = 3390
$ws.Range("L109").Value = 2832
$ws.Range("M109").Value = -2350
$ws.Range("N109").Value = -4912
$ws.Range("H120").Value = 17741.818
$ws.Range("I120").Value = 8201.666999999999
$ws.Range("K120").Value = 24605.001
$ws.Range("M120").Value = -19767.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 2288.75
$ws.Range("I80").Value = 2259.2307
$ws.Range("K80").Value = 2259.2307
$ws.Range("M80").Value = -1261.2307
$ws.Range("H83").Value = 2288.75
$ws.Range("I83").Value = 2259.2307
$ws.Range("K83").Value = 11296.1535
$ws.Range("M83").Value = -6304.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 82443.69500000001
$ws.Range("I7").Value = 104026.4
$ws.Range("J7").Value = 10501.333
$ws.Range("K7").Value = 104026.4
$ws.Range("L7").Value = 10501.333
$ws.Range("M7").Value = -103914.4
$ws.Range("N7").Value = -10725.333
$ws.Range("H22").Value = 1170.8077
$ws.Range("I22").Value = 1033.0714
$ws.Range("J22").Value = 1331.5
$ws.Range("K22").Value = 1033.0714
$ws.Range("L22").Value = 1331.5
$ws.Range("M22").Value = -738.0714
$ws.Range("N22").Value = -1921.5
$ws.Range("H27").Value = 1170.8077
$ws.Range("I27").Value = 1033.0714
$ws.Range("J27").Value = 1331.5
$ws.Range("K27").Value = 1033.0714
$ws.Range("L27").Value = 1331.5
$ws.Range("M27").Value = -926.0714
$ws.Range("N27").Value = -1545.5
$ws.Range("H40").Value = 7010.077
$ws.Range("I40").Value = 4663.6665
$ws.Range("K40").Value = 4663.6665
$ws.Range("M40").Value = -4527.6665
$ws.Range("H68").Value = 5662.0835
$ws.Range("I68").Value = 4085.2307
$ws.Range("J68").Value = 7525.636
$ws.Range("K68").Value = 4085.2307
$ws.Range("L68").Value = 7525.636
$ws.Range("M68").Value = -3336.2307
$ws.Range("N68").Value = -9023.636
$ws.Range("H71").Value = 5662.0835
$ws.Range("I71").Value = 4085.2307
$ws.Range("J71").Value = 7525.636
$ws.Range("K71").Value = 20426.1535
$ws.Range("L71").Value = 37628.18
$ws.Range("M71").Value = -16682.1535
$ws.Range("N71").Value = -45116.18
$ws.Range("H82").Value = 2110.6924
$ws.Range("I82").Value = 1269.1428
$ws.Range("K82").Value = 1269.1428
$ws.Range("M82").Value = -908.1428000000001
$ws.Range("H85").Value = 2110.6924
$ws.Range("I85").Value = 1269.1428
$ws.Range("K85").Value = 1269.1428
$ws.Range("M85").Value = -21.14280000000008
$ws.Range("H93").Value = 1710.6
$ws.Range("I93").Value = 1163.5264
$ws.Range("J93").Value = 3443
$ws.Range("K93").Value = 1163.5264
$ws.Range("L93").Value = 3443
$ws.Range("M93").Value = 84.47360000000003
$ws.Range("N93").Value = -5939
$ws.Range("H126").Value = 82443.69500000001
$ws.Range("I126").Value = 104026.4
$ws.Range("J126").Value = 10501.333
$ws.Range("K126").Value = 312079.2
$ws.Range("L126").Value = 31503.999
$ws.Range("M126").Value = -309609.2
$ws.Range("N126").Value = -36443.999
$ws.Range("H132").Value = 9007.432000000001
$ws.Range("I132").Value = 8214.967000000001
$ws.Range("J132").Value = 10705.571
$ws.Range("K132").Value = 24644.901
$ws.Range("L132").Value = 32116.713
$ws.Range("M132").Value = -22114.901
$ws.Range("N132").Value = -37176.713

Write-Output "Applied Moogle_Profits market data refresh."